# feat: add 2022-Q4 data
#
# The workbook tracks one "总计" (summary) sheet plus one sheet per
# quarter. This change adds a new "2022-Q4" quarter sheet (with its own
# fund-holding table) positioned right after "总计" and before the
# existing "2022-Q3" sheet, and appends a corresponding summary row on
# "总计".
#
# Implementation note: rather than inserting a brand new blank sheet for
# 2022-Q4 (which would mint a fresh, highest sheetId and leave it out of
# position), we duplicate the existing "2022-Q3" sheet so the duplicate
# keeps carrying the Q3 data forward (placed immediately after, and it
# picks up the new sheetId) while the original sheet slot - same sheetId,
# same position - is renamed to "2022-Q4" and its contents replaced with
# the new quarter's figures. This reproduces the intended sheetId/order
# layout: 总计=1, 2022-Q4=2, 2022-Q3=3.

$wb = $excel.ActiveWorkbook

$summary = $wb.Worksheets.Item(1)     # "总计"
$quarter = $wb.Worksheets.Item(2)     # currently "2022-Q3"

# 1. Duplicate the current "2022-Q3" sheet right after itself, so its data
#    survives unchanged in the new third slot.
$quarter.Copy($null, $quarter)
$quarterCopy = $wb.Worksheets.Item(3)

# 2. Rename: the original slot becomes the new quarter, the copy keeps the
#    old quarter's name. (Rename the original first so the copy can take
#    the now-freed "2022-Q3" name without a collision.)
$quarter.Name = "2022-Q4"
$quarterCopy.Name = "2022-Q3"

# 3. Overwrite the (now-named) "2022-Q4" sheet with the new quarter's fund
#    table. Copy/paste the header-style formatting from "总计" (which
#    already carries the right cell style) instead of re-deriving it.
$summary.Range("B1").Copy()
$quarter.Range("B1:H1").PasteSpecial(-4122)   # xlPasteFormats
$summary.Range("A2").Copy()
$quarter.Range("A2").PasteSpecial(-4122)      # xlPasteFormats

$quarter.Range("B1").Value = "基金代码"
$quarter.Range("C1").Value = "基金名称"
$quarter.Range("D1").Value = "基金规模"
$quarter.Range("E1").Value = "股票总仓位"
$quarter.Range("F1").Value = "仓位占比"
$quarter.Range("G1").Value = "持有市值(亿元)"
$quarter.Range("H1").Value = "仓位排名"

$quarter.Range("A2").Value = 0

# B2 and D2:G2 hold numeric-looking figures that must stay plain text (no
# implicit number coercion), same as the rest of the workbook's per-row
# data cells.
$quarter.Range("B2").NumberFormat = "@"
$quarter.Range("B2").Value = "159628"
$quarter.Range("C2").Value = "万家国证2000ETF"
$quarter.Range("D2:G2").NumberFormat = "@"
$quarter.Range("D2").Value = "2.55"
$quarter.Range("E2").Value = "97.28"
$quarter.Range("F2").Value = "0.51"
$quarter.Range("G2").Value = "0.0130"
$quarter.Range("H2").Value = 1

# Match "总计"'s page margins (inches -> points, x72) on the new sheet.
$quarter.PageSetup.LeftMargin = 54
$quarter.PageSetup.RightMargin = 54
$quarter.PageSetup.TopMargin = 72
$quarter.PageSetup.BottomMargin = 72
$quarter.PageSetup.HeaderMargin = 36
$quarter.PageSetup.FooterMargin = 36

# 4. Record the new quarter on "总计": relabel the existing row and append
#    a fresh row (with the same cell style as the row above) for the data
#    point that used to occupy it.
$summary.Range("B2").Value = "2022-Q4"

$summary.Range("A2").Copy()
$summary.Range("A3").PasteSpecial(-4122)      # xlPasteFormats
$summary.Range("A3").Value = 1
$summary.Range("B3").Value = "2022-Q3"
$summary.Range("C3").Value = 1
$summary.Range("D3").Value = 0.01
